$wb = $excel.ActiveWorkbook

# Sheet 1: AR
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.008517853992011859
$ws.Range("B3").Value = 0.7574397510546335
$ws.Range("B4").Value = 0.1449643810992358

# Sheet 2: SETAR
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B2").Value = -0.1924922941771267
$ws.Range("B3").Value = 0.5489540365738861
$ws.Range("B4").Value = 0.1056995749862125
$ws.Range("B5").Value = 0.1989364521800432
$ws.Range("B6").Value = 0.5554828186427432
$ws.Range("B7").Value = 0.1022394987945887

# Sheet 3: GARCH
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0008216282642889191
$ws.Range("B3").Value = 0.125904208037552
$ws.Range("B4").Value = 0.1036620094463338
$ws.Range("B5").Value = 0.1210785535428306

# Sheet 4: TARCH
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0006789516040753972
$ws.Range("B3").Value = 0.1275789815885949
$ws.Range("B4").Value = 0.08307317486530481
$ws.Range("B5").Value = 0.03876898703901042
$ws.Range("B6").Value = 0.1114728437573284

# Sheet 5: AR-TARCH
$ws = $wb.Worksheets.Item("AR-TARCH")
$ws.Range("B2").Value = 0.00935930655650029
$ws.Range("B3").Value = 0.7573059024682156
$ws.Range("B4").Value = 0.1289823215569525
$ws.Range("B5").Value = 0.07188003838808588
$ws.Range("B6").Value = 0.05942438202419818
$ws.Range("B7").Value = 0
